# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (same per-fund layout as the other
#    quarterly sheets) right before the "总计" (totals) sheet.
# 2. Insert a new top data row into "总计" for 2022-Q1 (2 funds, 0.91亿元),
#    pushing the existing history rows down and re-numbering the index
#    column.
#
# NOTE: worksheet object handles returned by this host resolve by
# *position*, not by a stable identity - e.g. after
# `$x = $wb.Worksheets.Add($total)`, the old `$total` handle now points
# at whatever sheet occupies that same slot (i.e. the newly-inserted
# one), not the original "总计" sheet. To stay safe, every sheet handle
# used below is (re)fetched by name immediately before it's used,
# rather than being cached across a structural operation (Add/Delete).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Scratch sheet / helper: Excel's Range.Value setter auto-coerces
# numeric-looking strings ("18.91", "009190", ...) to real numbers. The
# source data stores them as text, so stage every such string through a
# Text-formatted scratch cell and PasteSpecial *values only* into the
# destination - that keeps the destination cell's own number format
# (General) untouched while still landing a literal text value.
# ---------------------------------------------------------------------
$scratch = $wb.Worksheets.Add()
$scratch.Name = "__scratch__"
$scratchCell = $scratch.Cells.Item(1, 1)
$scratchCell.NumberFormat = "@"

function Set-TextValue($range, $text) {
    $scratchCell.Value = $text
    $scratchCell.Copy()
    $range.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# 1. New "2022-Q1" worksheet, inserted just before "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# Reuse the header/data formatting (fonts, borders, alignment) from an
# existing quarterly sheet with the same 8-column schema.
$srcQuarter = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Item("2022-Q1")
$srcQuarter.Range("A1:H3").Copy()
$q1.Range("A1:H3").PasteSpecial(-4122)

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

$q1.Cells.Item(2, 1).Value = 0
Set-TextValue $q1.Cells.Item(2, 2) "009190"
Set-TextValue $q1.Cells.Item(2, 3) "景顺长城核心优选一年持有期混合"
Set-TextValue $q1.Cells.Item(2, 4) "18.91"
Set-TextValue $q1.Cells.Item(2, 5) "89.60"
Set-TextValue $q1.Cells.Item(2, 6) "4.46"
Set-TextValue $q1.Cells.Item(2, 7) "0.8434"
$q1.Cells.Item(2, 8).Value = 6

$q1.Cells.Item(3, 1).Value = 1
Set-TextValue $q1.Cells.Item(3, 2) "008107"
Set-TextValue $q1.Cells.Item(3, 3) "华商医药医疗行业股票"
Set-TextValue $q1.Cells.Item(3, 4) "1.57"
Set-TextValue $q1.Cells.Item(3, 5) "91.35"
Set-TextValue $q1.Cells.Item(3, 6) "4.00"
Set-TextValue $q1.Cells.Item(3, 7) "0.0628"
$q1.Cells.Item(3, 8).Value = 4

# ---------------------------------------------------------------------
# 2. "总计" - insert a new top row for 2022-Q1, push history down
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total = $wb.Worksheets.Item("总计")
# Carry the row-2 formatting (index-column style) down from the row
# that just got pushed to row 3.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.91

# Re-number the index column (col A) for the rows that shifted down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
$total.Cells.Item(7, 1).Value = 5

# ---------------------------------------------------------------------
# Clean up the scratch sheet
# ---------------------------------------------------------------------
$scratch = $wb.Worksheets.Item("__scratch__")
$scratch.Delete() | Out-Null

Write-Output "done"
